$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: name / ID ---
$ws.Range("D1").Value = $null
$ws.Range("C1").Value = "Peter Parker"
$ws.Range("E1").NumberFormat = "@"
$ws.Range("E1").Value = "93027423234"
$ws.Range("E1").ClearFormats()

# --- Fall 2022 course list (col C/D), rows 6-8 ---
$ws.Range("C6").Value = "CYBR 3106"
$ws.Range("D6").Value = 3

$ws.Range("A7").Value = $null
$ws.Range("B7").Value = $null
$ws.Range("C7").Value = "CPSC 1302"
$ws.Range("D7").Value = 3

$ws.Range("C8").Value = $null
$ws.Range("D8").Value = $null

# --- Fall 2023 course list (rows 13-17) ---
$ws.Range("A13").Value = "CYBR 3108"
$ws.Range("C13").Value = "DSCI 3111"

$ws.Range("A14").Value = "CPSC 4115"
$ws.Range("C14").Value = "CPSC 3165"

$ws.Range("A15").Value = "CPSC 4155"
$ws.Range("C15").Value = "CPSC 4111"

$ws.Range("A16").Value = "CPSC 4157"
$ws.Range("C16").Value = "CPSC 4135"

$ws.Range("A17").Value = $null
$ws.Range("B17").Value = $null

# --- Fall 2024 course list (rows 22-24) ---
$ws.Range("A22").Value = "CPSC 4175"
$ws.Range("B22").Value = 3
$ws.Range("C22").Value = "CPSC 6985"
$ws.Range("D22").Value = 4

$ws.Range("A23").Value = "CPSC 6180"
$ws.Range("B23").Value = 3
$ws.Range("C23").Value = "CPSC 4000"
$ws.Range("D23").Value = 0

$ws.Range("A24").Value = "CPSC 6185"
$ws.Range("B24").Value = 3
